$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("M2").Value = 0.0006136666666666667
$ws.Range("N2").Value = 0.001841
$ws.Range("O2").Value = 0.000129696697123199
$ws.Range("P2").Value = 0.000129696697123199
$ws.Range("Q2").Value = 0.0262359735851111
$ws.Range("R2").Value = 0.236123762266
$ws.Range("S2").Value = 0.0001106220230120131
$ws.Range("T2").Value = 0.0001106220230120131

$ws.Range("G3").Value = 42.75280866666666
$ws.Range("H3").Value = 128.258426
$ws.Range("I3").Value = 0.8529286054750734
$ws.Range("J3").Value = 0.8529286054750735
$ws.Range("O3").Value = 0.8077214410831794
$ws.Range("P3").Value = 0.8077214410831794
$ws.Range("Q3").Value = 163.3916580948586
$ws.Range("R3").Value = 1470.524922853728
$ws.Range("S3").Value = 0.6889287223553929
$ws.Range("T3").Value = 0.688928722355393

$ws.Range("G4").Value = 42.75280866666666
$ws.Range("H4").Value = 128.258426
$ws.Range("I4").Value = 0.8529286054750734
$ws.Range("J4").Value = 0.8529286054750735
$ws.Range("O4").Value = 0.1921488622196973
$ws.Range("P4").Value = 0.1921488622196973
$ws.Range("Q4").Value = 38.86924328394021
$ws.Range("R4").Value = 349.823189555462
$ws.Range("S4").Value = 0.1638892610966685
$ws.Range("T4").Value = 0.1638892610966685

$ws.Range("I5").Value = 0.04642608686423023
$ws.Range("J5").Value = 0.04642608686423023
$ws.Range("M5").Value = 0.0006136666666666667
$ws.Range("N5").Value = 0.001841
$ws.Range("O5").Value = 0.000129696697123199
$ws.Range("P5").Value = 0.000129696697123199
$ws.Range("Q5").Value = 0.001428060427111111
$ws.Range("R5").Value = 0.012852543844
$ws.Range("S5").Value = 0.000006021310126645398
$ws.Range("T5").Value = 0.000006021310126645398

$ws.Range("I6").Value = 0.04642608686423023
$ws.Range("J6").Value = 0.04642608686423023
$ws.Range("O6").Value = 0.8077214410831794
$ws.Range("P6").Value = 0.8077214410831794
$ws.Range("Q6").Value = 8.893634546794667
$ws.Range("S6").Value = 0.03749934578582891
$ws.Range("T6").Value = 0.03749934578582891

$ws.Range("I7").Value = 0.04642608686423023
$ws.Range("J7").Value = 0.04642608686423023
$ws.Range("O7").Value = 0.1921488622196973
$ws.Range("P7").Value = 0.1921488622196973
$ws.Range("S7").Value = 0.008920719768274674
$ws.Range("T7").Value = 0.008920719768274674

$ws.Range("G8").Value = 5.044817999999999
$ws.Range("I8").Value = 0.1006453076606963
$ws.Range("J8").Value = 0.1006453076606963
$ws.Range("M8").Value = 0.0006136666666666667
$ws.Range("N8").Value = 0.001841
$ws.Range("O8").Value = 0.000129696697123199
$ws.Range("P8").Value = 0.000129696697123199
$ws.Range("Q8").Value = 0.003095836645999999
$ws.Range("R8").Value = 0.027862529814
$ws.Range("S8").Value = 0.00001305336398454051
$ws.Range("T8").Value = 0.00001305336398454051

$ws.Range("G9").Value = 5.044817999999999
$ws.Range("I9").Value = 0.1006453076606963
$ws.Range("J9").Value = 0.1006453076606963
$ws.Range("O9").Value = 0.8077214410831794
$ws.Range("P9").Value = 0.8077214410831794
$ws.Range("S9").Value = 0.08129337294195757
$ws.Range("T9").Value = 0.08129337294195757

$ws.Range("G10").Value = 5.044817999999999
$ws.Range("I10").Value = 0.1006453076606963
$ws.Range("J10").Value = 0.1006453076606963
$ws.Range("O10").Value = 0.1921488622196973
$ws.Range("P10").Value = 0.1921488622196973
$ws.Range("Q10").Value = 4.586558504121999
$ws.Range("R10").Value = 41.279026537098
$ws.Range("S10").Value = 0.01933888135475418
$ws.Range("T10").Value = 0.01933888135475418
